$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B5 value
$ws.Range("B5").Value = 56575

# Add new row 6 data
$ws.Range("A6").Value = 112446377
$ws.Range("B6").Value = 90152
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 1339
$ws.Range("F6").Value = "Brandticka"
$ws.Range("G6").Value = "Pycnoporellus fulgens"
$ws.Range("H6").Value = "(Fr.) Donk"
$ws.Range("I6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("P6").Value = "Molkomskogen, Vrm"
$ws.Range("Q6").Value = 427610
$ws.Range("R6").Value = 6608234
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Värmland"
$ws.Range("U6").Value = "Karlstad"
$ws.Range("V6").Value = "Värmland"
$ws.Range("W6").Value = "Nyed"

$ws.Range("Y6").Value = "'2023-10-01"
$ws.Range("Y6").Style = "Normal"
$ws.Range("AA6").Value = "'2023-10-01"
$ws.Range("AA6").Style = "Normal"

$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AT6").Value = ""
$ws.Range("AW6").Value = "Malin Max Nordgren"
$ws.Range("AX6").Value = "Malin Max Nordgren"
$ws.Range("AY6").Value = ""
